$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the daily series with 32 new rows (270-301), bringing the report
# up to date through 2021-06-28 ("aggiornamento fino al 28/06 incluso").
# Columns: A = date serial, B = nuovi pos., C = somma mobile 7gg.,
# D = somma mobile 7gg. per 100mila abitanti.
$colA = @(44344, 44345, 44346, 44347, 44348, 44349, 44350, 44351, 44352, 44353, 44354, 44355, 44356, 44357, 44358, 44359, 44360, 44361, 44362, 44363, 44364, 44365, 44366, 44367, 44368, 44369, 44370, 44371, 44372, 44373, 44374, 44375)
$colB = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 1, 0, 0, 0, 0, 1, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
$colC = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 1, 1, 1, 1, 1, 2, 2, 1, 1, 1, 1, 1, 0, 0, 0, 0, 0)
$colD = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 11.49954001839926, 11.49954001839926, 11.49954001839926, 11.49954001839926, 11.49954001839926, 22.99908003679853, 22.99908003679853, 11.49954001839926, 11.49954001839926, 11.49954001839926, 11.49954001839926, 11.49954001839926, 0, 0, 0, 0, 0)

$startRow = 270
$endRow = $startRow + $colA.Count - 1

# The date column (A) uses style index 2 (dedicated date number format,
# centered, thin border) in every existing row. Copy that formatting from
# the last pre-existing date cell (A269) onto the new date cells so the new
# rows are formatted exactly like the rows above them.
$ws.Range("A269").Copy()
$ws.Range("A270:A301").PasteSpecial(-4122)

for ($i = 0; $i -lt $colA.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $colA[$i]
    $ws.Cells.Item($r, 2).Value = $colB[$i]
    $ws.Cells.Item($r, 3).Value = $colC[$i]
    $ws.Cells.Item($r, 4).Value = $colD[$i]
}
